$d = $word.ActiveDocument

# --- Step 1: Fix typo "incon" -> "icon" and swap thread emoji for yarn emoji
#     in the "Put your R code..." paragraph ---
$d.Content.Find.Execute("the incon looks similar to 🧵", $true, $false, $false, $false, $false, $true, 1, $false, "the icon looks similar to 🧶", 2) | Out-Null

# --- Step 2: Move the "The icon resembles...filename." sentence (plus the
#     trailing space) from the start of the next paragraph onto the end of
#     the "Put your R code..." paragraph ---
$moveText = "The icon resembles a needle and thread as Knit implies. It is placed with other menu icons at the top, just under the filename."
$srcRange = $d.Content
$srcRange.Find.Execute($moveText) | Out-Null
$srcRange.MoveEnd(1, 1) | Out-Null   # include the trailing space
$srcRange.Delete()

$pR = $d.Paragraphs(6).Range
$pR.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
$pR.InsertAfter(" ")
$pR.Collapse(0) | Out-Null
$pR.InsertAfter($moveText)

# --- Step 3: Update the keyboard shortcut text ---
$d.Content.Find.Execute("Ctrl+Shift+P (or Cmd+Shift+P in Mac)", $true, $false, $false, $false, $false, $true, 1, $false, "Ctrl+Shift+K (or Cmd+Shift+ in Mac)", 2) | Out-Null

# --- Step 4: Remove ", then type Knit, before pressing Enter. " and split
#     the remainder off into its own new paragraph ---
$splitRange = $d.Content
$splitRange.Find.Execute(", then type Knit, before pressing Enter. ") | Out-Null
$splitRange.Delete()
$breakPoint = $d.Range($splitRange.Start, $splitRange.Start)
$breakPoint.InsertParagraphAfter()
